# Add a "time_taken" column (F) with per-row timestamps to the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell, re-using the bold/bordered header style (same as B1:E1).
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "time_taken"

$values = @(
    "2021-10-05 13:41:07.460240",
    "2021-10-05 13:41:07.460253",
    "2021-10-05 13:41:07.460257",
    "2021-10-05 13:41:07.460260",
    "2021-10-05 13:41:07.460264",
    "2021-10-05 13:41:07.460267",
    "2021-10-05 13:41:07.460270",
    "2021-10-05 13:41:07.460273",
    "2021-10-05 13:41:07.460276",
    "2021-10-05 13:41:07.460280",
    "2021-10-05 13:41:07.460283",
    "2021-10-05 13:41:07.460286",
    "2021-10-05 13:41:07.460289",
    "2021-10-05 13:41:07.460292",
    "2021-10-05 13:41:07.460295",
    "2021-10-05 13:41:07.460298",
    "2021-10-05 13:41:07.460301",
    "2021-10-05 13:41:07.460305",
    "2021-10-05 13:41:07.460308",
    "2021-10-05 13:41:07.460311",
    "2021-10-05 13:41:07.460314",
    "2021-10-05 13:41:07.460317",
    "2021-10-05 13:41:07.460320",
    "2021-10-05 13:41:07.460323",
    "2021-10-05 13:41:07.460326",
    "2021-10-05 13:41:07.460330",
    "2021-10-05 13:41:07.460333",
    "2021-10-05 13:41:07.460336",
    "2021-10-05 13:41:07.460339",
    "2021-10-05 13:41:07.460342",
    "2021-10-05 13:41:07.460345",
    "2021-10-05 13:41:07.460348",
    "2021-10-05 13:41:07.460352",
    "2021-10-05 13:41:07.460355",
    "2021-10-05 13:41:07.460358",
    "2021-10-05 13:41:07.460361",
    "2021-10-05 13:41:07.460364",
    "2021-10-05 13:41:07.460367",
    "2021-10-05 13:41:07.460370",
    "2021-10-05 13:41:07.460373",
    "2021-10-05 13:41:07.460377",
    "2021-10-05 13:41:07.460380",
    "2021-10-05 13:41:07.460383",
    "2021-10-05 13:41:07.460387"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $values[$i]
}
